$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.911.62"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.636.27"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.868.40"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.639.63"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "27.932.28"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "1.402.45"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.560"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.853"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "1.776.41"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.06%  "
